# Append four more "moses/bro/1234/m@g.c/Male/0" rows (52-55) to the Users
# sheet, matching the existing rows 6-51 exactly (values + shared-string
# typing + default styling). Copy/paste from the last existing data row so
# the pasted cells keep the same text-type (not auto-coerced to numbers)
# and the same (default) cell style as the source.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

$lastRow = 51
$firstNewRow = 52
$lastNewRow = 55

$ws.Range("A$lastRow`:F$lastRow").Copy() | Out-Null

for ($r = $firstNewRow; $r -le $lastNewRow; $r++) {
    $ws.Range("A$r`:F$r").PasteSpecial() | Out-Null
}

$excel.CutCopyMode = 0
